# Applies the diff:
#   1. "Description (50 words max): "   -> "Description :"          (merge 2 runs into 1, keeps bold rPr)
#   2. "Generated Summary (100 words max):" -> "Generated Summary :" (keeps the colon/space split point)
#   3. Splits the run "<EMDASH>where users compare their lives negatively with peers<EMDASH>emerged as..."
#      into four runs with identical formatting:
#         "<EMDASH>" / "where users compare their lives negatively with peers" / "<EMDASH>" / "emerged as..."
#   4. Splits the run "...social media addiction<EMDASH>unhealthy comparisons<EMDASH>could reduce..."
#      into five runs with identical formatting:
#         "...addiction" / "<EMDASH>" / "unhealthy comparisons" / "<EMDASH>" / "could reduce..."

$d = $word.ActiveDocument
$emdash = [string][char]0x2014

# --- helper: force a run split at a Range's boundaries without changing its
#     visible formatting. The COM shim only splits a run on a genuine property
#     write (a same-value no-op gets optimized away and the run stays merged
#     with its identically-formatted neighbours), so bounce the font size to a
#     different value and back to force the split while leaving the original
#     formatting intact.
function Split-AtRange($rng) {
    $origSize = $rng.Font.Size
    if ($origSize -eq 13) {
        $rng.Font.Size = 14
    } else {
        $rng.Font.Size = 13
    }
    $rng.Font.Size = $origSize
}

# 1. "Description (50 words max): " -> "Description :"
$d.Content.Find.Execute("Description (50 words max): ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Description :", 2) | Out-Null

# 2. "Generated Summary (100 words max):" -> "Generated Summary :"
#    (the following run's leading ": " keeps only its space -> " The research...")
$d.Content.Find.Execute("Generated Summary (100 words max):", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Generated Summary :", 2) | Out-Null

# 3. Split the em-dash-delimited clause in the "Unfavorable social comparisons..." paragraph
$anchor3 = $d.Content
$anchor3.Find.Execute("Unfavorable social comparisons", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

$dash3a = $d.Range($anchor3.End, $d.Content.End)
$dash3a.Find.Execute($emdash, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-AtRange $dash3a

$mid3 = $d.Range($dash3a.End, $d.Content.End)
$mid3.Find.Execute("where users compare their lives negatively with peers", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
Split-AtRange $mid3

$dash3b = $d.Range($mid3.End, $d.Content.End)
$dash3b.Find.Execute($emdash, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-AtRange $dash3b

# 4. Split the em-dash-delimited clause in the "...social media addiction...comparisons..." paragraph
$anchor4 = $d.Content
$anchor4.Find.Execute("social media addiction", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

$dash4a = $d.Range($anchor4.End, $d.Content.End)
$dash4a.Find.Execute($emdash, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-AtRange $dash4a

$mid4 = $d.Range($dash4a.End, $d.Content.End)
$mid4.Find.Execute("unhealthy comparisons", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
Split-AtRange $mid4

$dash4b = $d.Range($mid4.End, $d.Content.End)
$dash4b.Find.Execute($emdash, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-AtRange $dash4b

Write-Host "Done"
